$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 369.660016
$ws.Range("H2").Value = 1108.980048
$ws.Range("I2").Value = 0.1642435133179984
$ws.Range("J2").Value = 0.1642435133179984
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 23.63579766666667
$ws.Range("N2").Value = 70.907393
$ws.Range("O2").Value = 0.06827844587621175
$ws.Range("P2").Value = 0.06827844587621175
$ws.Range("Q2").Value = 8737.209343632763
$ws.Range("R2").Value = 78634.88409269488
$ws.Range("S2").Value = 0.01121429183460182
$ws.Range("T2").Value = 0.01121429183460182

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 369.660016
$ws.Range("H3").Value = 1108.980048
$ws.Range("I3").Value = 0.1642435133179984
$ws.Range("J3").Value = 0.1642435133179984
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 181.2883913333334
$ws.Range("N3").Value = 543.865174
$ws.Range("O3").Value = 0.5237009467675041
$ws.Range("P3").Value = 0.523700946767504
$ws.Range("Q3").Value = 67015.06964089427
$ws.Range("R3").Value = 603135.6267680485
$ws.Range("S3").Value = 0.08601448342505692
$ws.Range("T3").Value = 0.08601448342505692

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 369.660016
$ws.Range("H4").Value = 1108.980048
$ws.Range("I4").Value = 0.1642435133179984
$ws.Range("J4").Value = 0.1642435133179984
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 111.1005463333333
$ws.Range("N4").Value = 333.301639
$ws.Range("O4").Value = 0.3209442197221123
$ws.Range("P4").Value = 0.3209442197221123
$ws.Range("Q4").Value = 41069.42973518874
$ws.Range("R4").Value = 369624.8676166987
$ws.Range("S4").Value = 0.05271300622626335
$ws.Range("T4").Value = 0.05271300622626336

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 369.660016
$ws.Range("H5").Value = 1108.980048
$ws.Range("I5").Value = 0.1642435133179984
$ws.Range("J5").Value = 0.1642435133179984
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 30.14303933333333
$ws.Range("N5").Value = 90.42911799999999
$ws.Range("O5").Value = 0.08707638763417187
$ws.Range("P5").Value = 0.08707638763417187
$ws.Range("Q5").Value = 11142.67640224863
$ws.Range("R5").Value = 100284.0876202377
$ws.Range("S5").Value = 0.0143017318320763
$ws.Range("T5").Value = 0.0143017318320763

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 249.925644
$ws.Range("H6").Value = 749.776932
$ws.Range("I6").Value = 0.1110443760810293
$ws.Range("J6").Value = 0.1110443760810294
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 23.63579766666667
$ws.Range("N6").Value = 70.907393
$ws.Range("O6").Value = 0.06827844587621175
$ws.Range("P6").Value = 0.06827844587621175
$ws.Range("Q6").Value = 5907.191953295363
$ws.Range("R6").Value = 53164.72757965828
$ws.Range("S6").Value = 0.007581937422106264
$ws.Range("T6").Value = 0.007581937422106265

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 249.925644
$ws.Range("H7").Value = 749.776932
$ws.Range("I7").Value = 0.1110443760810293
$ws.Range("J7").Value = 0.1110443760810294
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 181.2883913333334
$ws.Range("N7").Value = 543.865174
$ws.Range("O7").Value = 0.5237009467675041
$ws.Range("P7").Value = 0.523700946767504
$ws.Range("Q7").Value = 45308.61795370735
$ws.Range("R7").Value = 407777.5615833662
$ws.Range("S7").Value = 0.05815404488684185
$ws.Range("T7").Value = 0.05815404488684184

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 249.925644
$ws.Range("H8").Value = 749.776932
$ws.Range("I8").Value = 0.1110443760810293
$ws.Range("J8").Value = 0.1110443760810294
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 111.1005463333333
$ws.Range("N8").Value = 333.301639
$ws.Range("O8").Value = 0.3209442197221123
$ws.Range("P8").Value = 0.3209442197221123
$ws.Range("Q8").Value = 27766.87559111017
$ws.Range("R8").Value = 249901.8803199915
$ws.Range("S8").Value = 0.03563905063585475
$ws.Range("T8").Value = 0.03563905063585476

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 249.925644
$ws.Range("H9").Value = 749.776932
$ws.Range("I9").Value = 0.1110443760810293
$ws.Range("J9").Value = 0.1110443760810294
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 30.14303933333333
$ws.Range("N9").Value = 90.42911799999999
$ws.Range("O9").Value = 0.08707638763417187
$ws.Range("P9").Value = 0.08707638763417187
$ws.Range("Q9").Value = 7533.518517500664
$ws.Range("R9").Value = 67801.66665750596
$ws.Range("S9").Value = 0.009669343136226473
$ws.Range("T9").Value = 0.009669343136226475

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1552.995524333333
$ws.Range("H10").Value = 4658.986573
$ws.Range("I10").Value = 0.6900109020274287
$ws.Range("J10").Value = 0.6900109020274288
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 23.63579766666667
$ws.Range("N10").Value = 70.907393
$ws.Range("O10").Value = 0.06827844587621175
$ws.Range("P10").Value = 0.06827844587621175
$ws.Range("Q10").Value = 36706.28799038158
$ws.Range("R10").Value = 330356.5919134342
$ws.Range("S10").Value = 0.04711287202807583
$ws.Range("T10").Value = 0.04711287202807584

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1552.995524333333
$ws.Range("H11").Value = 4658.986573
$ws.Range("I11").Value = 0.6900109020274287
$ws.Range("J11").Value = 0.6900109020274288
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 181.2883913333334
$ws.Range("N11").Value = 543.865174
$ws.Range("O11").Value = 0.5237009467675041
$ws.Range("P11").Value = 0.523700946767504
$ws.Range("Q11").Value = 281540.0603542565
$ws.Range("R11").Value = 2533860.543188309
$ws.Range("S11").Value = 0.3613593626716639
$ws.Range("T11").Value = 0.3613593626716639

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1552.995524333333
$ws.Range("H12").Value = 4658.986573
$ws.Range("I12").Value = 0.6900109020274287
$ws.Range("J12").Value = 0.6900109020274288
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 111.1005463333333
$ws.Range("N12").Value = 333.301639
$ws.Range("O12").Value = 0.3209442197221123
$ws.Range("P12").Value = 0.3209442197221123
$ws.Range("Q12").Value = 172538.6512066548
$ws.Range("R12").Value = 1552847.860859893
$ws.Range("S12").Value = 0.221455010550944
$ws.Range("T12").Value = 0.221455010550944

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1552.995524333333
$ws.Range("H13").Value = 4658.986573
$ws.Range("I13").Value = 0.6900109020274287
$ws.Range("J13").Value = 0.6900109020274288
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 30.14303933333333
$ws.Range("N13").Value = 90.42911799999999
$ws.Range("O13").Value = 0.08707638763417187
$ws.Range("P13").Value = 0.08707638763417187
$ws.Range("Q13").Value = 46812.00517447029
$ws.Range("R13").Value = 421308.0465702325
$ws.Range("S13").Value = 0.06008365677674497
$ws.Range("T13").Value = 0.06008365677674498

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 78.101406
$ws.Range("H14").Value = 234.304218
$ws.Range("I14").Value = 0.03470120857354342
$ws.Range("J14").Value = 0.03470120857354343
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 23.63579766666667
$ws.Range("N14").Value = 70.907393
$ws.Range("O14").Value = 0.06827844587621175
$ws.Range("P14").Value = 0.06827844587621175
$ws.Range("Q14").Value = 1845.989029698186
$ws.Range("R14").Value = 16613.90126728367
$ws.Range("S14").Value = 0.002369344591427819
$ws.Range("T14").Value = 0.00236934459142782

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 78.101406
$ws.Range("H15").Value = 234.304218
$ws.Range("I15").Value = 0.03470120857354342
$ws.Range("J15").Value = 0.03470120857354343
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 181.2883913333334
$ws.Range("N15").Value = 543.865174
$ws.Range("O15").Value = 0.5237009467675041
$ws.Range("P15").Value = 0.523700946767504
$ws.Range("Q15").Value = 14158.87825461155
$ws.Range("R15").Value = 127429.9042915039
$ws.Range("S15").Value = 0.01817305578394132
$ws.Range("T15").Value = 0.01817305578394132

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 78.101406
$ws.Range("H16").Value = 234.304218
$ws.Range("I16").Value = 0.03470120857354342
$ws.Range("J16").Value = 0.03470120857354343
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 111.1005463333333
$ws.Range("N16").Value = 333.301639
$ws.Range("O16").Value = 0.3209442197221123
$ws.Range("P16").Value = 0.3209442197221123
$ws.Range("Q16").Value = 8677.108876001477
$ws.Range("R16").Value = 78093.97988401329
$ws.Range("S16").Value = 0.01113715230905016
$ws.Range("T16").Value = 0.01113715230905017

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 78.101406
$ws.Range("H17").Value = 234.304218
$ws.Range("I17").Value = 0.03470120857354342
$ws.Range("J17").Value = 0.03470120857354343
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 30.14303933333333
$ws.Range("N17").Value = 90.42911799999999
$ws.Range("O17").Value = 0.08707638763417187
$ws.Range("P17").Value = 0.08707638763417187
$ws.Range("Q17").Value = 2354.213753046636
$ws.Range("R17").Value = 21187.92377741972
$ws.Range("S17").Value = 0.003021655889124115
$ws.Range("T17").Value = 0.003021655889124116
